$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 swap values for columns D, J, K, L, M, P
# (date, volumen, precio minimo, precio maximo, precio promedio ponderado, precio $/Kg)

# New row 2 values (previously held by row 3)
$ws.Range("D2").Value = 44827
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 31000
$ws.Range("M2").Value = 30500
$ws.Range("P2").Value = 1220

# New row 3 values (previously held by row 2)
$ws.Range("D3").Value = 44414
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 31000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31500
$ws.Range("P3").Value = 1260
